# Weekly fruit/vegetable price update: insert two new daily records
# (row 64 and 65) for "Ají" (Inferno variety) on 2022-05-30, pushing the
# existing rows 64-80 down to 66-82.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the current row 64 (formatting carries over
# from the row above, same as Excel's native Insert behaviour).
$ws.Rows.Item(64).Resize(2).Insert()

# ---- New row 64 ----
$ws.Range("A64").Value = 1
$ws.Range("B64").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C64").Value = "Arica y Parinacota"
$ws.Range("D64").Value = 44711
$ws.Range("E64").Value = 15
$ws.Range("F64").Value = 100112021
$ws.Range("G64").Value = "Ají"
$ws.Range("H64").Value = "Inferno"
$ws.Range("I64").Value = "Primera"
$ws.Range("J64").Value = 140
$ws.Range("K64").Value = 13000
$ws.Range("L64").Value = 14000
$ws.Range("M64").Value = 13500
$ws.Range("N64").Value = "$/caja 15 kilos"
$ws.Range("O64").Value = "Región de Arica y Parinacota"
$ws.Range("P64").Value = 900
$ws.Range("Q64").Value = 15
$ws.Range("R64").Value = "Hortaliza"

# ---- New row 65 ----
$ws.Range("A65").Value = 1
$ws.Range("B65").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C65").Value = "Arica y Parinacota"
$ws.Range("D65").Value = 44711
$ws.Range("E65").Value = 15
$ws.Range("F65").Value = 100112021
$ws.Range("G65").Value = "Ají"
$ws.Range("H65").Value = "Inferno"
$ws.Range("I65").Value = "Segunda"
$ws.Range("J65").Value = 120
$ws.Range("K65").Value = 12000
$ws.Range("L65").Value = 13000
$ws.Range("M65").Value = 12500
$ws.Range("N65").Value = "$/caja 15 kilos"
$ws.Range("O65").Value = "Región de Arica y Parinacota"
$ws.Range("P65").Value = 833
$ws.Range("Q65").Value = 15
$ws.Range("R65").Value = "Hortaliza"
